# Insert a new weekly price record as row 31 ("Hortaliza, Feria Lagunitas de
# Puerto Montt - Poroto verde"), pushing all existing records (old rows
# 31-102) down by one row (new rows 32-103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31:102 down to 32:103, leaving a blank row 31 to populate.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new record's data.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44838
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 100112031
$ws.Range("G31").Value = "Poroto verde"
$ws.Range("H31").Value = "Magnum"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 35
$ws.Range("K31").Value = 32000
$ws.Range("L31").Value = 32000
$ws.Range("M31").Value = 32000
$ws.Range("N31").Value = "`$/malla 25 kilos"
$ws.Range("O31").Value = "Perú"
$ws.Range("P31").Value = 1280
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
